$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Superficie" (B) and "Population" (C) columns were originally stored
# as shared-string text (cell style s="2" applies a Text "@" number format).
# The edit replaces those text labels with real numeric values while the
# cells keep their existing (Text) style. Because the Text format otherwise
# forces any newly-entered value back into a string, we briefly switch the
# number format to a non-text one, write the numbers, then restore "@".
$dataRange = $ws.Range("B2:C19")
$dataRange.NumberFormat = "0"

$ws.Range("B2").Value = 69711
$ws.Range("C2").Value = 7948287

$ws.Range("B3").Value = 47784
$ws.Range("C3").Value = 2811423

$ws.Range("B4").Value = 27208
$ws.Range("C4").Value = 3318904

$ws.Range("B5").Value = 39151
$ws.Range("C5").Value = 2576252

$ws.Range("B6").Value = 8680
$ws.Range("C6").Value = 334938

$ws.Range("B7").Value = 57433
$ws.Range("C7").Value = 5549586

$ws.Range("B8").Value = 1703
$ws.Range("C8").Value = 390253

$ws.Range("B9").Value = 83534
$ws.Range("C9").Value = 268700

$ws.Range("B10").Value = 31813
$ws.Range("C10").Value = 6003815

$ws.Range("B11").Value = 12011
$ws.Range("C11").Value = 12174880

$ws.Range("B12").Value = 1128
$ws.Range("C12").Value = 372594

$ws.Range("B13").Value = 2504
$ws.Range("C13").Value = 250143

$ws.Range("B14").Value = 29906
$ws.Range("C14").Value = 3330478

$ws.Range("B15").Value = 83809
$ws.Range("C15").Value = 5956978

$ws.Range("B16").Value = 72724
$ws.Range("C16").Value = 5845102

$ws.Range("B17").Value = 32082
$ws.Range("C17").Value = 3757600

$ws.Range("B18").Value = 2505
$ws.Range("C18").Value = 853659

$ws.Range("B19").Value = 31400
$ws.Range("C19").Value = 5030890

# Restore the original Text number format on that range.
$dataRange.NumberFormat = "@"

# Match the new selection left behind in the workbook (B2:C19, the table
# of values that was just pasted in).
$ws.Range("B2:C19").Select() | Out-Null
